# Weekly update: insert two new rows of data (week of 44511) above the
# existing row 125, pushing all the old rows (125-158) down to (127-160).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above row 125; Excel copies the formatting
# (e.g. the date style on column D) from the row above automatically.
$ws.Rows("125:126").Insert()

# New row 125: Terminal Hortofrutícola Agro Chillán, Pera, Primera, 44511
$ws.Range("A125").Value = 7
$ws.Range("B125").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C125").Value = "Ñuble"
$ws.Range("D125").Value = 44511
$ws.Range("E125").Value = 16
$ws.Range("F125").Value = "Fruta"
$ws.Range("G125").Value = 100104
$ws.Range("H125").Value = "Frutos de pepita"
$ws.Range("I125").Value = 100104005
$ws.Range("J125").Value = "Pera"
$ws.Range("K125").Value = "Packham's Triumph"
$ws.Range("L125").Value = "Primera"
$ws.Range("M125").Value = 120
$ws.Range("N125").Value = 9500
$ws.Range("O125").Value = 10000
$ws.Range("P125").Value = 9750
$ws.Range("Q125").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R125").Value = "Provincia de Curicó"
$ws.Range("S125").Value = 609
$ws.Range("T125").Value = 16

# New row 126: Terminal Hortofrutícola Agro Chillán, Pera, Segunda, 44511
$ws.Range("A126").Value = 7
$ws.Range("B126").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C126").Value = "Ñuble"
$ws.Range("D126").Value = 44511
$ws.Range("E126").Value = 16
$ws.Range("F126").Value = "Fruta"
$ws.Range("G126").Value = 100104
$ws.Range("H126").Value = "Frutos de pepita"
$ws.Range("I126").Value = 100104005
$ws.Range("J126").Value = "Pera"
$ws.Range("K126").Value = "Packham's Triumph"
$ws.Range("L126").Value = "Segunda"
$ws.Range("M126").Value = 120
$ws.Range("N126").Value = 8000
$ws.Range("O126").Value = 8500
$ws.Range("P126").Value = 8250
$ws.Range("Q126").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R126").Value = "Provincia de Curicó"
$ws.Range("S126").Value = 516
$ws.Range("T126").Value = 16
